$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose values look like plain numbers
# (so Excel stores them as literal text, matching the source data which
# uses locale-formatted strings such as "0.830" / "2.20" / "8.10" rather
# than numeric values) before writing the new value.

$ws.Range('D2').Value = '44.178.49'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '2.242.22'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.33'
$ws.Range('E5').Value = '  -2.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.66'
$ws.Range('E6').Value = '  -3.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.66'
$ws.Range('E10').Value = '  -3.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0809'
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.19'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').Value = '2.582.16'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '2.328.37'
$ws.Range('E15').Value = '  +3.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.830'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.52'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '43.915.57'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = '0.0₃0963'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.38'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.07'
$ws.Range('E21').Value = '  -7.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.59'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.09'
$ws.Range('E23').Value = '  +4.04%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '238.53'
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.99'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.92'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').Value = '  +3.18%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.07'
$ws.Range('E29').Value = '  +4.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.05'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.86'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.04'
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0794'
$ws.Range('E33').Value = '  -4.19%  '
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.19'
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('E38').Value = '  -7.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.56'
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.82'
$ws.Range('E40').Value = '  -4.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.31'
$ws.Range('E41').Value = '  -7.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0299'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('D44').Value = '1.738.76'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.56'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.191'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '14.99'
$ws.Range('E47').Value = '  +8.22%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '99.87'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.92'
$ws.Range('E49').Value = '  -4.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.10'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.57'
$ws.Range('E51').Value = '  -2.24%  '
